function Set-TextValue {
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '34.683.27'
$ws.Range("E2").Value = '  +2.00%  '
# Row 3
$ws.Range("D3").Value = '1.794.43'
$ws.Range("E3").Value = '  +0.48%  '
# Row 4
$ws.Range("E4").Value = '  +0.22%  '
# Row 5
Set-TextValue $ws 5 4 '225.55'
$ws.Range("E5").Value = '  +0.24%  '
# Row 6
Set-TextValue $ws 6 4 '0.554'
$ws.Range("E6").Value = '  -0.49%  '
# Row 7
$ws.Range("E7").Value = '  +0.19%  '
# Row 8
Set-TextValue $ws 8 4 '32.79'
$ws.Range("E8").Value = '  +6.42%  '
# Row 9
$ws.Range("E9").Value = '  +2.10%  '
# Row 10
Set-TextValue $ws 10 4 '0.0669'
$ws.Range("E10").Value = '  +1.23%  '
# Row 11
$ws.Range("E11").Value = '  +1.21%  '
# Row 12
$ws.Range("D12").Value = '2.052.65'
$ws.Range("E12").Value = '  +0.48%  '
# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 13 4 '11.11'
$ws.Range("E13").Value = '  +11.14%  '
# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.803.21'
$ws.Range("E14").Value = '  +0.97%  '
# Row 15
$ws.Range("E15").Value = '  +1.32%  '
# Row 16
$ws.Range("D16").Value = '34.696.72'
$ws.Range("E16").Value = '  +2.22%  '
# Row 17
Set-TextValue $ws 17 4 '4.29'
$ws.Range("E17").Value = '  +2.14%  '
# Row 18
Set-TextValue $ws 18 4 '69.20'
$ws.Range("E18").Value = '  +1.01%  '
# Row 19
Set-TextValue $ws 19 4 '254.74'
$ws.Range("E19").Value = '  +1.18%  '
# Row 20
$ws.Range("D20").Value = '0.0₃0764'
$ws.Range("E20").Value = '  +3.31%  '
# Row 21
Set-TextValue $ws 21 4 '1.00'
$ws.Range("E21").Value = '  +0.30%  '
# Row 22
Set-TextValue $ws 22 4 '10.41'
$ws.Range("E22").Value = '  +1.21%  '
# Row 23
Set-TextValue $ws 23 4 '4.24'
$ws.Range("E23").Value = '  +0.42%  '
# Row 24
Set-TextValue $ws 24 4 '2.14'
$ws.Range("E24").Value = '  -0.96%  '
# Row 25
Set-TextValue $ws 25 4 '158.58'
$ws.Range("E25").Value = '  +0.20%  '
# Row 26
Set-TextValue $ws 26 4 '16.43'
$ws.Range("E26").Value = '  -0.22%  '
# Row 27
Set-TextValue $ws 27 4 '7.10'
$ws.Range("E27").Value = '  +2.26%  '
# Row 28
$ws.Range("E28").Value = '  -0.50%  '
# Row 29
$ws.Range("E29").Value = '  +0.11%  '
# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 30 4 '3.77'
$ws.Range("E30").Value = '  -0.85%  '
# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 31 4 '0.0519'
$ws.Range("E31").Value = '  +1.32%  '
# Row 32
$ws.Range("E32").Value = '  +0.25%  '
# Row 33
Set-TextValue $ws 33 4 '3.59'
$ws.Range("E33").Value = '  +1.71%  '
# Row 34
Set-TextValue $ws 34 4 '1.87'
$ws.Range("E34").Value = '  +6.54%  '
# Row 35
$ws.Range("D35").Value = '1.451.28'
# Row 36
$ws.Range("E36").Value = '  +0.50%  '
# Row 37
$ws.Range("E37").Value = '  +2.08%  '
# Row 38
Set-TextValue $ws 38 4 '0.629'
$ws.Range("E38").Value = '  +0.14%  '
# Row 39
Set-TextValue $ws 39 4 '83.20'
$ws.Range("E39").Value = '  +0.05%  '
# Row 40
Set-TextValue $ws 40 4 '2.83'
$ws.Range("E40").Value = '  +4.41%  '
# Row 41
$ws.Range("E41").Value = '  +0.07%  '
# Row 42
Set-TextValue $ws 42 4 '0.900'
$ws.Range("E42").Value = '  +1.57%  '
# Row 43
$ws.Range("E43").Value = '  -0.48%  '
# Row 44
$ws.Range("E44").Value = '  -0.63%  '
# Row 45
Set-TextValue $ws 45 4 '5.95'
$ws.Range("E45").Value = '  +3.65%  '
# Row 46
$ws.Range("E46").Value = '  -1.33%  '
# Row 47
$ws.Range("D47").Value = '1.951.92'
# Row 48
Set-TextValue $ws 48 4 '104.31'
$ws.Range("E48").Value = '  +6.26%  '
# Row 49
$ws.Range("E49").Value = '  +0.23%  '
# Row 50
Set-TextValue $ws 50 4 '11.93'
$ws.Range("E50").Value = '  -0.01%  '
# Row 51
$ws.Range("D51").Value = '0.0₆0123'
$ws.Range("E51").Value = '  +5.43%  '
